$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-15 Tuesday" "2025-04-22 Tuesday"

Replace-Text "262×5=1310" "733×6=4398"
Replace-Text "771×5=3855" "926×8=7408"
Replace-Text "614×6=3684" "556×9=5004"
Replace-Text "138×7=966" "635×4=2540"
Replace-Text "452×7=3164" "381×8=3048"
Replace-Text "806×2=1612" "853×6=5118"
Replace-Text "188×8=1504" "808×4=3232"
Replace-Text "394×7=2758" "179×7=1253"
Replace-Text "765×3=2295" "613×8=4904"
Replace-Text "832×7=5824" "999×9=8991"
Replace-Text "502×5=2510" "980×9=8820"
Replace-Text "553×2=1106" "561×4=2244"
Replace-Text "782×5=3910" "753×2=1506"
Replace-Text "718×3=2154" "257×3=771"
Replace-Text "923×8=7384" "643×7=4501"
Replace-Text "892×7=6244" "632×4=2528"
Replace-Text "548×7=3836" "711×8=5688"
Replace-Text "975×8=7800" "902×9=8118"
Replace-Text "929×8=7432" "352×5=1760"
Replace-Text "641×5=3205" "465×4=1860"
Replace-Text "325×6=1950" "544×6=3264"
Replace-Text "336×4=1344" "680×9=6120"
Replace-Text "776×7=5432" "307×7=2149"
Replace-Text "586×6=3516" "516×3=1548"
Replace-Text "489×5=2445" "602×3=1806"
